# Weekly update: insert a new price-observation row for the week of
# 2022-07-05 just above the existing row 225 (old rows 225-232 shift
# down to 226-233).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down by inserting a fresh row at 225.
$ws.Rows.Item(225).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(225, 1).Value = 4
$ws.Cells.Item(225, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(225, 3).Value = "Los Lagos"
$ws.Cells.Item(225, 4).Value = 44747
$ws.Cells.Item(225, 5).Value = 10
$ws.Cells.Item(225, 6).Value = 100112039
$ws.Cells.Item(225, 7).Value = "Ciboulette"
$ws.Cells.Item(225, 8).Value = "Sin especificar"
$ws.Cells.Item(225, 9).Value = "Primera"
$ws.Cells.Item(225, 10).Value = 240
$ws.Cells.Item(225, 11).Value = 2500
$ws.Cells.Item(225, 12).Value = 3000
$ws.Cells.Item(225, 13).Value = 2750
$ws.Cells.Item(225, 14).Value = "`$/docena de atados"
$ws.Cells.Item(225, 15).Value = "Región Metropolitana"
$ws.Cells.Item(225, 16).Value = 917
$ws.Cells.Item(225, 17).Value = 3
$ws.Cells.Item(225, 18).Value = "Hortaliza"
